$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing data
# (currently in columns A:O) one column to the right (B:P) and keeps
# their formatting/values intact.
$ws.Columns.Item(1).Insert()

# Add the new "ix" header in the now-empty column A, row 1 (no special
# style - matches the other un-styled header cell).
$ws.Cells.Item(1, 1).Value = "ix"

# Fill column A (rows 2-12) with a simple running index 1..11.
for ($i = 2; $i -le 12; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}

# Update the active selection on the sheet to match the new state.
$ws.Range("A17").Select()
